$d = $word.ActiveDocument

# Turn on revision tracking for the duration of the edits so the newly
# touched text is kept in its own run instead of silently re-merging with
# the neighbouring run (which shares identical formatting). We accept each
# resulting revision individually afterwards (rather than AcceptAllRevisions,
# which also touches/normalizes unrelated parts of the document), so the
# saved XML ends up with plain split runs and no tracked-change markup.
$d.TrackRevisions = $true

# --- Fix 1: "Dit houd in dat" -> "Dit houdt in dat" --------------------
# Insert a "t" right after "Dit houd"; the rest of the sentence before and
# after the insertion point is left untouched.
$r1 = $d.Content
$r1.Find.Execute("Dit houd", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ins1 = $r1.Duplicate
$ins1.Collapse(0)
$ins1.InsertAfter("t")

# --- Fix 2: "aan het eind word er gekeken" -> "aan het eind wordt er gekeken"
# Replace just the standalone word "word" with "wordt".
$r2 = $d.Content
$r2.Find.Execute("eind word er", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wStart = $r2.Start + 5   # skip past "eind "
$wEnd = $wStart + 4       # length of "word"
$wordRange = $d.Range($wStart, $wEnd)
$wordRange.Text = "wordt"

$d.TrackRevisions = $false
foreach ($rev in $d.Revisions) {
    $rev.Accept()
}
